# Update "想去人数" (want-to-go count) figures on both the 展览 and
# 全部类型 worksheets, as produced by the latest gh-pages data refresh.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 11614
    $ws.Range("F20").Value = 467
    $ws.Range("F22").Value = 10995
}
